# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) holds a recomputed strike-count value (s_vals)
# derived from a fresh std/mean calculation. The new values below are the
# result of that recalculation for each data row (rows 2-68).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strike count) values, one per data row, starting at row 2.
$sVals = @(0,1,2,1,0,2,1,1,0,1,0,0,0,3,1,1,1,2,2,0,0,1,0,1,0,2,0,1,1,1,1,1,0,0,1,3,1,0,2,0,2,1,1,0,0,0,2,1,2,0,0,1,1,2,2,0,0,2,2,1,2,1,1,2,1,1,3)

$startRow = 2
for ($i = 0; $i -lt $sVals.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $sVals[$i]
}
